{"js": "// Resume update for Brooks job application.\nconst body = context.document.body;\n\n// 1) Bio sentence: drop \"14 years of\", add \"and analysis\" before \", and design for manufacturing.\"\nlet bioIntro = body.search(\n  \"Mechanical Engineer with 14 years of industry experience in structural design and design for manufacturing.\",\n  { matchCase: true }\n);\nbioIntro.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < bioIntro.items.length; i++) {\n  bioIntro.items[i].insertText(\n    \"Mechanical Engineer with industry experience in structural design and analysis, and design for manufacturing.\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\n// 2) Master's thesis line: \"nanomaterials\" -> \"porous nanomaterials\"\nlet thesisLine = body.search(\n  \"Master\\u2019s Degree in Mechanical Engineering, with thesis work on numerical modeling of nanomaterials.\",\n  { matchCase: true }\n);\nthesisLine.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < thesisLine.items.length; i++) {\n  thesisLine.items[i].insertText(\n    \"Master\\u2019s Degree in Mechanical Engineering, with thesis work on numerical modeling of porous nanomaterials.\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\n// 3) Remove the \"Dedicated to working toward the clean energy transition.\" sentence,\n//    along with the single leading space that separated it from the previous sentence.\nlet dedicatedSentence = body.search(\n  \" Dedicated to working toward the clean energy transition.\",\n  { matchCase: true }\n);\ndedicatedSentence.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < dedicatedSentence.items.length; i++) {\n  dedicatedSentence.items[i].delete();\n}\nawait context.sync();\n\n// 4) Fix typo \"deflectino\" -> \"deflection\"\nlet deflectionTypo = body.search(\"blade deflectino and test optimal structural parameters.\", {\n  matchCase: true,\n});\ndeflectionTypo.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < deflectionTypo.items.length; i++) {\n  deflectionTypo.items[i].insertText(\n    \"blade deflection and test optimal structural parameters.\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\n// 5) \"Created model, drawing, ...\" -> \"Created CAD models, drawings, ...\"\nlet catiaLine = body.search(\n  \"Created model, drawing, and CNC templates for steel wall panels using Catia 3DExperience.\",\n  { matchCase: true }\n);\ncatiaLine.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < catiaLine.items.length; i++) {\n  catiaLine.items[i].insertText(\n    \"Created CAD models, drawings, and CNC templates for steel wall panels using Catia 3DExperience.\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n", "ps1": "# Resume update for Brooks job application.\n$d = $word.ActiveDocument\n\n# 1) Bio sentence: drop \"14 years of\", add \"and analysis\" before \", and design for manufacturing.\"\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$r.Find.Text = \"Mechanical Engineer with 14 years of industry experience in structural design and design for manufacturing.\"\n$r.Find.Replacement.Text = \"Mechanical Engineer with industry experience in structural design and analysis, and design for manufacturing.\"\n$r.Find.Execute([ref]\"\", [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]\"\", [ref]2)\n\n# 2) Master's thesis line: \"nanomaterials\" -> \"porous nanomaterials\"\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$r.Find.Text = \"Master\u2019s Degree in Mechanical Engineering, with thesis work on numerical modeling of nanomaterials.\"\n$r.Find.Replacement.Text = \"Master\u2019s Degree in Mechanical Engineering, with thesis work on numerical modeling of porous nanomaterials.\"\n$r.Find.Execute([ref]\"\", [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]\"\", [ref]2)\n\n# 3) Remove the \"Dedicated to working toward the clean energy transition.\" sentence,\n#    along with the single leading space that separated it from the previous sentence.\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Text = \" Dedicated to working toward the clean energy transition.\"\n$r.Find.Execute()\nif ($r.Find.Found) {\n    $r.Delete()\n}\n\n# 4) Fix typo \"deflectino\" -> \"deflection\"\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$r.Find.Text = \"blade deflectino and test optimal structural parameters.\"\n$r.Find.Replacement.Text = \"blade deflection and test optimal structural parameters.\"\n$r.Find.Execute([ref]\"\", [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]\"\", [ref]2)\n\n# 5) \"Created model, drawing, ...\" -> \"Created CAD models, drawings, ...\"\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Replacement.ClearFormatting()\n$r.Find.Text = \"Created model, drawing, and CNC templates for steel wall panels using Catia 3DExperience.\"\n$r.Find.Replacement.Text = \"Created CAD models, drawings, and CNC templates for steel wall panels using Catia 3DExperience.\"\n$r.Find.Execute([ref]\"\", [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]\"\", [ref]2)\n"}
